$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 20,10
$arr[0,0] = -16.87794310902494
$arr[0,1] = -0.2490846644997531
$arr[0,2] = -16.87794310902494
$arr[0,3] = -16.87794310902494
$arr[0,4] = -16.87794310902494
$arr[0,5] = -16.87794310902494
$arr[0,6] = -16.87794310902494
$arr[0,7] = -16.87794310902494
$arr[0,8] = -16.87794310902494
$arr[0,9] = -16.87794310902494
$arr[1,0] = -16.87794310902494
$arr[1,1] = -16.87794310902494
$arr[1,2] = -16.87794310902494
$arr[1,3] = -16.87794310902494
$arr[1,4] = -16.87794310902494
$arr[1,5] = -16.87794310902494
$arr[1,6] = -16.87794310902494
$arr[1,7] = 0.6273179651234523
$arr[1,8] = -16.87794310902494
$arr[1,9] = -16.87794310902494
$arr[2,0] = -16.87794310902494
$arr[2,1] = -0.309556451500452
$arr[2,2] = 0.2103171875544621
$arr[2,3] = -16.87794310902494
$arr[2,4] = 3.939263075407696
$arr[2,5] = -16.87794310902494
$arr[2,6] = 1.324147226562418
$arr[2,7] = -16.87794310902494
$arr[2,8] = 2.089815077060935
$arr[2,9] = -16.87794310902494
$arr[3,0] = -16.87794310902494
$arr[3,1] = 0.3630438232280809
$arr[3,2] = -16.87794310902494
$arr[3,3] = -16.87794310902494
$arr[3,4] = -16.87794310902494
$arr[3,5] = 3.444338263762611
$arr[3,6] = -16.87794310902494
$arr[3,7] = -16.87794310902494
$arr[3,8] = -16.87794310902494
$arr[3,9] = -16.87794310902494
$arr[4,0] = -16.87794310902494
$arr[4,1] = -16.87794310902494
$arr[4,2] = -16.87794310902494
$arr[4,3] = -16.87794310902494
$arr[4,4] = -16.87794310902494
$arr[4,5] = -16.87794310902494
$arr[4,6] = -16.87794310902494
$arr[4,7] = -16.87794310902494
$arr[4,8] = -16.87794310902494
$arr[4,9] = -16.87794310902494
$arr[5,0] = 3.134939179006455
$arr[5,1] = -16.87794310902494
$arr[5,2] = -16.87794310902494
$arr[5,3] = -16.87794310902494
$arr[5,4] = -16.87794310902494
$arr[5,5] = -16.87794310902494
$arr[5,6] = -16.87794310902494
$arr[5,7] = -16.87794310902494
$arr[5,8] = -16.87794310902494
$arr[5,9] = -16.87794310902494
$arr[6,0] = -16.87794310902494
$arr[6,1] = -16.87794310902494
$arr[6,2] = -16.87794310902494
$arr[6,3] = 1.358374789950199
$arr[6,4] = -16.87794310902494
$arr[6,5] = -16.87794310902494
$arr[6,6] = -16.87794310902494
$arr[6,7] = -16.87794310902494
$arr[6,8] = -16.87794310902494
$arr[6,9] = -16.87794310902494
$arr[7,0] = 3.4874192520429
$arr[7,1] = -16.87794310902494
$arr[7,2] = -16.87794310902494
$arr[7,3] = -16.87794310902494
$arr[7,4] = -16.87794310902494
$arr[7,5] = -16.87794310902494
$arr[7,6] = -16.87794310902494
$arr[7,7] = -16.87794310902494
$arr[7,8] = -16.87794310902494
$arr[7,9] = -16.87794310902494
$arr[8,0] = -16.87794310902494
$arr[8,1] = -16.87794310902494
$arr[8,2] = -16.87794310902494
$arr[8,3] = -16.87794310902494
$arr[8,4] = -16.87794310902494
$arr[8,5] = -16.87794310902494
$arr[8,6] = -16.87794310902494
$arr[8,7] = 1.294274395010331
$arr[8,8] = -16.87794310902494
$arr[8,9] = -16.87794310902494
$arr[9,0] = -16.87794310902494
$arr[9,1] = -16.87794310902494
$arr[9,2] = -16.87794310902494
$arr[9,3] = 2.098959386035207
$arr[9,4] = -16.87794310902494
$arr[9,5] = 1.841997126120907
$arr[9,6] = -16.87794310902494
$arr[9,7] = -16.87794310902494
$arr[9,8] = -16.87794310902494
$arr[9,9] = 4.321916715137496
$arr[10,0] = -16.87794310902494
$arr[10,1] = -16.87794310902494
$arr[10,2] = -16.87794310902494
$arr[10,3] = -16.87794310902494
$arr[10,4] = -16.87794310902494
$arr[10,5] = -16.87794310902494
$arr[10,6] = -16.87794310902494
$arr[10,7] = -16.87794310902494
$arr[10,8] = -16.87794310902494
$arr[10,9] = -16.87794310902494
$arr[11,0] = -16.87794310902494
$arr[11,1] = -16.87794310902494
$arr[11,2] = -16.87794310902494
$arr[11,3] = 1.750263889381972
$arr[11,4] = -16.87794310902494
$arr[11,5] = -16.87794310902494
$arr[11,6] = -16.87794310902494
$arr[11,7] = -16.87794310902494
$arr[11,8] = 1.512347623181374
$arr[11,9] = -16.87794310902494
$arr[12,0] = -16.87794310902494
$arr[12,1] = -16.87794310902494
$arr[12,2] = 0.7964491823379056
$arr[12,3] = -16.87794310902494
$arr[12,4] = -16.87794310902494
$arr[12,5] = -16.87794310902494
$arr[12,6] = -16.87794310902494
$arr[12,7] = -16.87794310902494
$arr[12,8] = -16.87794310902494
$arr[12,9] = -16.87794310902494
$arr[13,0] = -16.87794310902494
$arr[13,1] = -16.87794310902494
$arr[13,2] = -0.03244080698583635
$arr[13,3] = -16.87794310902494
$arr[13,4] = -16.87794310902494
$arr[13,5] = -16.87794310902494
$arr[13,6] = -16.87794310902494
$arr[13,7] = -16.87794310902494
$arr[13,8] = -16.87794310902494
$arr[13,9] = -16.87794310902494
$arr[14,0] = -16.87794310902494
$arr[14,1] = -16.87794310902494
$arr[14,2] = -16.87794310902494
$arr[14,3] = -16.87794310902494
$arr[14,4] = -16.87794310902494
$arr[14,5] = -16.87794310902494
$arr[14,6] = -16.87794310902494
$arr[14,7] = -16.87794310902494
$arr[14,8] = 2.095207695351244
$arr[14,9] = -16.87794310902494
$arr[15,0] = -16.87794310902494
$arr[15,1] = 0.3108871701000157
$arr[15,2] = -0.13973207268387
$arr[15,3] = -16.87794310902494
$arr[15,4] = -16.87794310902494
$arr[15,5] = -16.87794310902494
$arr[15,6] = 1.957098787804147
$arr[15,7] = 0.2892888073571042
$arr[15,8] = 2.450269420607891
$arr[15,9] = -16.87794310902494
$arr[16,0] = -16.87794310902494
$arr[16,1] = -16.87794310902494
$arr[16,2] = -16.87794310902494
$arr[16,3] = -16.87794310902494
$arr[16,4] = -16.87794310902494
$arr[16,5] = -16.87794310902494
$arr[16,6] = 2.376847937423416
$arr[16,7] = -0.003814244182196341
$arr[16,8] = 1.656333172582617
$arr[16,9] = -16.87794310902494
$arr[17,0] = -16.87794310902494
$arr[17,1] = -16.87794310902494
$arr[17,2] = 3.007808367521533
$arr[17,3] = -16.87794310902494
$arr[17,4] = -16.87794310902494
$arr[17,5] = -16.87794310902494
$arr[17,6] = 1.951596924315694
$arr[17,7] = 0.9498027411744683
$arr[17,8] = -16.87794310902494
$arr[17,9] = -16.87794310902494
$arr[18,0] = -16.87794310902494
$arr[18,1] = 3.215243599496967
$arr[18,2] = 2.843431486461498
$arr[18,3] = -16.87794310902494
$arr[18,4] = 2.220163529626624
$arr[18,5] = -16.87794310902494
$arr[18,6] = 1.238812034807572
$arr[18,7] = 3.567016148611036
$arr[18,8] = -16.87794310902494
$arr[18,9] = -16.87794310902494
$arr[19,0] = -16.87794310902494
$arr[19,1] = 2.708795910134028
$arr[19,2] = -16.87794310902494
$arr[19,3] = 3.290970126862925
$arr[19,4] = -16.87794310902494
$arr[19,5] = 2.467091730729547
$arr[19,6] = 1.131694182890404
$arr[19,7] = -16.87794310902494
$arr[19,8] = -16.87794310902494
$arr[19,9] = -16.87794310902494
$ws.Range("B2:K21").Value = $arr
